$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E) so periods run in ascending order
# (2103 .. 2109) from row 16 to row 22 - part of refreshing the EC database.
$ws.Range("E16").Value = "2103"
$ws.Range("E17").Value = "2104"
$ws.Range("E18").Value = "2105"
$ws.Range("E19").Value = "2106"
$ws.Range("E20").Value = "2107"
$ws.Range("E21").Value = "2108"
$ws.Range("E22").Value = "2109"

# Swap the "Valor Mora" (F) values between the first and last rows to match
# the new ordering of periods.
$ws.Range("F16").Value = 36341
$ws.Range("F22").Value = 24227
